$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "T973C, C1009T, G1006T, G942A, A949G, C958T, A1045T, G960A, G1043A, C952T"
$ws.Range("B3").Value  = "C8407T, C8480T, T8473C, A8456G, G8488A, C8481T, C8476T, C8502T, A8413G"
$ws.Range("B4").Value  = "G12191T, C12295T, G12202T, C12312T, C12171T, C12325T, C12242T, G12223T, C12213T, G12223A, G12188T, G12163T, A12173G, C12241T, G12163A, G12243A"
$ws.Range("B5").Value  = "T12628C, T12640C, C12559T, G12565T, A12644G, C12624A, A12544G, C12616T, A12634G, C12623T, G12569A"
$ws.Range("B6").Value  = "C13424T, C13366T, C13426T, C13432T, C13372T, C13384T, C13436T, C13452T, C13378T, C13381T, C13424A, C13446T, C13423T, C13405T, C13421T"
$ws.Range("B7").Value  = "G16795T, A16797G, C16767T, C16806T, C16762T, G16827A, G16846T, T16800C, C16751A"
$ws.Range("B8").Value  = "T17865C, A17858G, G17944T, C17934T, C17862T, A17853T, T17963C, C17883T, G17893T, C17913T, A17916G"
$ws.Range("B9").Value  = "C19512T, C19488T, C19402T, G19518T, A19476G, G19419A, A19515G, T19519C, C19495T, T19521C, G19509T, G19509A, G19480T, T19452A, T19404C"
$ws.Range("B10").Value = "T24224C, G24257T, C24319T, A24253T, C24237T, G24328T, C24223T, G24236T, G24332T"
$ws.Range("B11").Value = "C24700T, G24697A, A24673G, C24707A, T24733C, T24739C, C24715T, T24703C"
$ws.Range("B12").Value = "C26408T, A26319G, C26340T, A26358G, G26389T, C26366T, G26314A, G26325A"
